# Update the Fixed-Residue PSSM matrix (B2:K21) with the supplemental-figure values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pssm = @{}
$pssm[2] = @(-17.50505943745325, -0.6664340836241436, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[3] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -0.04581194581039263, -17.50505943745325, -17.50505943745325)
$pssm[4] = @(-17.50505943745325, -0.3962389416708416, 0.05796548667358671, -17.50505943745325, 4.013394424981032, -17.50505943745325, 4.321920726859038, -17.50505943745325, 3.630850835131421, -17.50505943745325)
$pssm[5] = @(-17.50505943745325, 0.04283671358832233, -17.50505943745325, -17.50505943745325, -17.50505943745325, 3.879716409391777, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[6] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[7] = @(3.481955396159227, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[8] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, 1.185037459326044, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[9] = @(3.14189384628868, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[10] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -0.9075742373891454, -17.50505943745325, 1.451909363654354)
$pssm[11] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, 2.885543036478572, -17.50505943745325, 0.9796572705965259, -17.50505943745325, -17.50505943745325, -17.50505943745325, 0.6970308981341095)
$pssm[12] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[13] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, 1.544100456330312, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, 0.382571733528961, 3.639128894133241)
$pssm[14] = @(-17.50505943745325, -17.50505943745325, 0.6500601902443305, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, 0.8186330831425837)
$pssm[15] = @(-17.50505943745325, -17.50505943745325, -0.5757894688061733, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)
$pssm[16] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, 0.9844475293578412, -17.50505943745325)
$pssm[17] = @(-17.50505943745325, 0.05841243899554328, -0.4244812965237328, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -0.8614709600723072, 1.125069014932662, -17.50505943745325)
$pssm[18] = @(-17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -0.5180064260760033, 1.103502219182092, -17.50505943745325)
$pssm[19] = @(-17.50505943745325, -17.50505943745325, 2.640510805409841, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325, -0.4918131060428424, -17.50505943745325, -17.50505943745325)
$pssm[20] = @(-17.50505943745325, 3.474173846889748, 3.283646948761439, -17.50505943745325, 1.945106710753218, -17.50505943745325, -17.50505943745325, 4.047733940226188, -17.50505943745325, 0.5061077481123389)
$pssm[21] = @(-17.50505943745325, 2.439779853380195, -17.50505943745325, 2.891450230615131, -17.50505943745325, 1.725873374222416, -17.50505943745325, -17.50505943745325, -17.50505943745325, -17.50505943745325)

foreach ($row in $pssm.Keys) {
    $vals = $pssm[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}
